$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "('Astral Drift', ['{2}{W}', 'Enchantment', 'Whenever you cycle Astral Drift or cycle another card while Astral Drift is on the battlefield, you may exile target creature. If you do, return that card to the battlefield under its owner" + [char]0x2019 + "s control at the beginning of the next end step.', 'Cycling {2}{W} ({2}{W}, Discard this card: Draw a card.)'])"

$ws.Range("A2").Value = $newText

# Remove rows 3 to 6
$ws.Range("A3:A6").EntireRow.Delete()
